# Apply weekly refresh of Fruta/Hortaliza data: rows are reshuffled to
# reflect a new set of reported dates/values (rows keep their other
# columns, but D (Fecha), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado) and S (Precio $/Kg)
# are updated per row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values keyed by target row number: Fecha, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, PrecioKg
$rows = @{
    2  = @(44517, "Especial", 100, 27000, 27000, 27000, 2700)
    3  = @(44517, "Primera",  30,  25000, 25000, 25000, 2500)
    4  = @(44432, "Primera",  20,  20000, 20000, 20000, 2000)
    5  = @(44473, "Primera",  180, 20000, 20000, 20000, 2000)
    6  = @(44476, "Primera",  120, 20000, 20000, 20000, 2000)
    7  = @(44466, "Primera",  60,  20000, 20000, 20000, 2000)
    8  = @(44434, "Primera",  20,  20000, 20000, 20000, 2000)
    9  = @(44511, "Primera",  120, 28000, 28000, 28000, 2800)
    10 = @(44435, "Primera",  40,  20000, 20000, 20000, 2000)
    11 = @(44503, "Primera",  60,  30000, 30000, 30000, 3000)
    12 = @(44503, "Segunda",  50,  25000, 25000, 25000, 2500)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]

    $ws.Cells.Item($r, 4).Value  = $vals[0]  # D - Fecha
    $ws.Cells.Item($r, 12).Value = $vals[1]  # L - Calidad
    $ws.Cells.Item($r, 13).Value = $vals[2]  # M - Volumen
    $ws.Cells.Item($r, 14).Value = $vals[3]  # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $vals[4]  # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $vals[5]  # P - Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $vals[6]  # S - Precio $/Kg
}
